# Update NATMI cell-cell communication values per "Natmi following Dr Hou advice".
# Columns A-D (Sending cluster, Ligand symbol, Receptor symbol, Target cluster) are unchanged.
# Columns E-T for rows 2-13 are updated with recomputed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$data = @{
    2 = @(3, 1, 212.4076363333333, 637.222909, 0.5080632835800084, 0.5080632835800084, 3, 1, 0.3116673333333334, 0.935002, 0.0414413620607491, 0.0414413620607491, 66.20052159564645, 595.804694360818, 0.02105483448461217, 0.02105483448461217)
    3 = @(3, 1, 212.4076363333333, 637.222909, 0.5080632835800084, 0.5080632835800084, 3, 1, 3.794584666666667, 11.383754, 0.5045532214096876, 0.5045532214096876, 805.9987599133763, 7253.988839220386, 0.2563449664102769, 0.2563449664102769)
    4 = @(3, 1, 212.4076363333333, 637.222909, 0.5080632835800084, 0.5080632835800084, 3, 1, 3.414430666666667, 10.243292, 0.4540054165295632, 0.4540054165295633, 725.2511473307142, 6527.260325976427, 0.2306634826851193, 0.2306634826851194)
    5 = @(3, 1, 161.9384456666667, 485.815337, 0.3873447295187379, 0.3873447295187379, 3, 1, 0.3116673333333334, 0.935002, 0.0414413620607491, 0.0414413620607491, 50.47092352507489, 454.238311725674, 0.01605209317830895, 0.01605209317830895)
    6 = @(3, 1, 161.9384456666667, 485.815337, 0.3873447295187379, 0.3873447295187379, 3, 1, 3.794584666666667, 11.383754, 0.5045532214096876, 0.5045532214096876, 614.4891428705664, 5530.402285835098, 0.1954360310747433, 0.1954360310747433)
    7 = @(3, 1, 161.9384456666667, 485.815337, 0.3873447295187379, 0.3873447295187379, 3, 1, 3.414430666666667, 10.243292, 0.4540054165295632, 0.4540054165295633, 552.9275949966004, 4976.348354969404, 0.1758566052656856, 0.1758566052656856)
    8 = @(3, 1, 0.338549, 1.015647, 0.0008097840527861261, 0.0008097840527861261, 3, 1, 0.3116673333333334, 0.935002, 0.0414413620607491, 0.0414413620607491, 0.1055146640326667, 0.949631976294, 0.00003355855412253061, 0.00003355855412253061)
    9 = @(3, 1, 0.338549, 1.015647, 0.0008097840527861261, 0.0008097840527861261, 3, 1, 3.794584666666667, 11.383754, 0.5045532214096876, 0.5045532214096876, 1.284652844315333, 11.561875598838, 0.0004085791524794325, 0.0004085791524794325)
    10 = @(3, 1, 0.338549, 1.015647, 0.0008097840527861261, 0.0008097840527861261, 3, 1, 3.414430666666667, 10.243292, 0.4540054165295632, 0.4540054165295633, 1.155952087769333, 10.403568789924, 0.000367646346184163, 0.0003676463461841631)
    11 = @(3, 1, 43.38855633333333, 130.165669, 0.1037822028484675, 0.1037822028484675, 3, 1, 0.3116673333333334, 0.935002, 0.0414413620607491, 0.0414413620607491, 13.52279564959311, 121.705160846338, 0.004300875843705446, 0.004300875843705446)
    12 = @(3, 1, 43.38855633333333, 130.165669, 0.1037822028484675, 0.1037822028484675, 3, 1, 3.794584666666667, 11.383754, 0.5045532214096876, 0.5045532214096876, 164.6415505712696, 1481.773955141426, 0.05236364477218791, 0.05236364477218791)
    13 = @(3, 1, 43.38855633333333, 130.165669, 0.1037822028484675, 0.1037822028484675, 3, 1, 3.414430666666667, 10.243292, 0.4540054165295632, 0.4540054165295633, 148.1472173269275, 1333.324955942348, 0.04711768223257409, 0.04711768223257409)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}
